$wb = $excel.ActiveWorkbook

# --- Update the daily conversion message on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.68 = 26669.25 pesos`n✅ 26669.25 pesos = 6.65 = 957.7 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate values on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 149.798
$wsTasas.Range("O10").Value = 3995
$wsTasas.Range("N12").Value = 4009.99
$wsTasas.Range("O12").Value = 144
